$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Block 1 (Fall/Spring/Summer 2022) - header row 3 stays the same.
# The course list in the data rows (4-7) is replaced with a new set of
# courses, and the old row 8 (CYBR 3119) is removed entirely.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "CPSC 1301K"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "CPSC 1302"
$ws.Range("D4").Value = 3

$ws.Range("A5").Value = "CPSC 4111"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "CYBR 2159"
$ws.Range("D5").Value = 3

$ws.Range("A6").Value = "CPSC 6180"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "CYBR 3115"
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = "CPSC 6185"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 6985"
$ws.Range("D7").Value = 4

# Row 8 (previously held "CYBR 3119" / 3 credits) is cleared so it drops out
# of the sheet entirely, leaving row 11 as the Total row for this block.
$ws.Range("A8:F8").ClearContents()

# ---------------------------------------------------------------------------
# Block 2 (Fall/Spring/Summer 2023) - brand new block of rows.
# Header row 12, data rows 13-16, totals row 20.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Fall 2023"
$ws.Range("B12").Value = "Credits"
$ws.Range("C12").Value = "Spring 2023"
$ws.Range("D12").Value = "Credits"
$ws.Range("E12").Value = "Summer 2023"
$ws.Range("F12").Value = "Credits"

$ws.Range("A13").Value = "CPSC 2108"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CYBR 3106"
$ws.Range("D13").Value = 3

$ws.Range("A14").Value = "CYBR 2160"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "CYBR 3108"
$ws.Range("D14").Value = 3

$ws.Range("C15").Value = "CYBR 3119"
$ws.Range("D15").Value = 3

$ws.Range("C16").Value = "CPSC 4115"
$ws.Range("D16").Value = 3

$ws.Range("A20").Value = "Total"
$ws.Range("B20").Formula = "=SUM(B13:B19)"
$ws.Range("C20").Value = "Total"
$ws.Range("D20").Formula = "=SUM(D13:D19)"
$ws.Range("E20").Value = "Total"
$ws.Range("F20").Formula = "=SUM(F13:F19)"

# ---------------------------------------------------------------------------
# Block 3 (Fall/Spring/Summer 2024) - header only, no course rows yet.
# Header row 21, totals row 29.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Fall 2024"
$ws.Range("B21").Value = "Credits"
$ws.Range("C21").Value = "Spring 2024"
$ws.Range("D21").Value = "Credits"
$ws.Range("E21").Value = "Summer 2024"
$ws.Range("F21").Value = "Credits"

$ws.Range("A29").Value = "Total"
$ws.Range("B29").Formula = "=SUM(B22:B28)"
$ws.Range("C29").Value = "Total"
$ws.Range("D29").Formula = "=SUM(D22:D28)"
$ws.Range("E29").Value = "Total"
$ws.Range("F29").Formula = "=SUM(F22:F28)"

# ---------------------------------------------------------------------------
# Block 4 (Fall/Spring/Summer 2025) - header only, no course rows yet.
# Header row 30, totals row 38.
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = "Fall 2025"
$ws.Range("B30").Value = "Credits"
$ws.Range("C30").Value = "Spring 2025"
$ws.Range("D30").Value = "Credits"
$ws.Range("E30").Value = "Summer 2025"
$ws.Range("F30").Value = "Credits"

$ws.Range("A38").Value = "Total"
$ws.Range("B38").Formula = "=SUM(B31:B37)"
$ws.Range("C38").Value = "Total"
$ws.Range("D38").Formula = "=SUM(D31:D37)"
$ws.Range("E38").Value = "Total"
$ws.Range("F38").Formula = "=SUM(F31:F37)"
